# Add placeno to calculations
#
# Turns the lone "{placeno}" run into a docxtemplater loop:
#   {#calculations}{placeno}{/}
# split across runs exactly like the sibling "{#calculations}...{/}"
# blocks already present in this template (separate <w:r> per literal
# chunk, all sharing the surrounding formatting).

$d = $word.ActiveDocument

# Locate the lone "{placeno}" placeholder.
$target = $d.Content
$target.Find.ClearFormatting()
$target.Find.Execute("{placeno}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$placenoStart = $target.Start
$placenoEnd = $target.End

# --- 1) Append "{/}" right after "{placeno}" -----------------------------
# Inserting directly at the end of the existing run merges into it,
# inheriting its formatting (incl. lang="en-US"), matching the diff's
# "{placeno}{/}" run.
$closeTag = $d.Range($placenoEnd, $placenoEnd)
$closeTag.InsertBefore("{/}")

# --- 2) Insert "{#calculations}" right before "{placeno}" ----------------
# The run immediately preceding "{placeno}" is an existing empty run
# (no text); inserting into its collapsed start reuses that run's
# formatting (no lang attribute), matching the diff's new runs.
$prefixStart = $placenoStart
$prefixRange = $d.Range($prefixStart, $prefixStart)
$prefixRange.InsertBefore("{#calculations}")
$prefixEnd = $prefixStart + 15   # Len("{#calculations}")

# --- 3) Split the single "{#calculations}" run into separate runs --------
# "{#" | "calculations" | "}"   (matching the diff's run boundaries)
# A momentary formatting round-trip (Underline off then back on) forces
# the engine to keep the two sides as distinct runs instead of
# re-coalescing them, while leaving the final formatting unchanged.
$splitAfterHash = $prefixStart + 2
$afterHash = $d.Range($splitAfterHash, $prefixEnd)
$afterHash.Font.Underline = 0
$afterHash.Font.Underline = 1

$splitBeforeBrace = $prefixEnd - 1
$beforeBrace = $d.Range($splitBeforeBrace, $prefixEnd)
$beforeBrace.Font.Underline = 0
$beforeBrace.Font.Underline = 1
